$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: strip trailing "$" from period labels (rows 2-33) ---
$periods = @("Q4 2018", "Q1 2019", "Q2 2019", "Q3 2019", "Q4 2019", "Q1 2020", "Q2 2020", "Q3 2020", "Q4 2020", "Q1 2021", "Q2 2021", "Q3 2021", "Q4 2021", "Q1 2022", "Q2 2022", "Q3 2022", "Q4 2022", "Q1 2023", "Q2 2023", "Q3 2023", "Q4 2023", "Q1 2024", "Q2 2024", "Q3 2024", "Q4 2024", "Q1 2025", "Q2 2025", "Q3 2025", "Q4 2025", "Q1 2026", "Q2 2026", "Q3 2026")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $periods[$i]
}

# --- Column B: updated gw_simul values (rows 7-33; row 6 remains 0) ---
$bVals = @(0.1738722793249711, 0.5150643219552331, 0.5449564712044563, 0.6656359317161844, 0.6394326790202198, 0.5499030255340499, 0.4739436644373797, 0.481250667138948, 0.4861263726029907, 0.4617416225808446, 0.4429031399077412, 0.4300569271249406, 0.4199716787195968, 0.4094504198574394, 0.3993059629861925, 0.3920267936153891, 0.3866437324893752, 0.3815354709359833, 0.377077635881538, 0.3737445216481247, 0.3711149410618871, 0.368801882201169, 0.3668570920831674, 0.3653358537959684, 0.3641064195189391, 0.3630584604805778, 0.3621881907178166)
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $row = $i + 7
    $ws.Cells.Item($row, 2).Value = $bVals[$i]
}

# --- Column C: updated gcpi_simul values (rows 6-33) ---
$cVals = @(7.67906559498353, 4.202493577422769, 2.38438084022046, 1.428444879444643, 0.9945465253977196, 0.7903392530878948, 0.6631425388183607, 0.5850249005666273, 0.5492101456520386, 0.5220886387941499, 0.4955214444981157, 0.4754453870865816, 0.4591209121674681, 0.4444740153276895, 0.4313726624522503, 0.4194237022554943, 0.4092583715129271, 0.4008747531416141, 0.3935086634951455, 0.3872582702529465, 0.382184298024455, 0.3779704759810946, 0.3744596304291303, 0.3715618091135703, 0.369188936447792, 0.3672504804985987, 0.3656458770444848, 0.3643236168463379)
for ($i = 0; $i -lt $cVals.Length; $i++) {
    $row = $i + 6
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
}

# --- Column D: updated cf1_simul values (rows 6-33) ---
$dVals = @(0.5611116717935283, 1.489088372285891, 1.085193212321585, 0.713969258209547, 0.5490062612939669, 0.5641744746785298, 0.5244735419505069, 0.4510081718839272, 0.419497192063879, 0.4228322453919932, 0.4140146534186994, 0.3950978251922602, 0.3867267274047179, 0.3858679692419142, 0.3818286423547494, 0.3758359617018224, 0.3725651627338414, 0.3710718402870521, 0.3689620342107532, 0.3666308393821051, 0.3651331201992724, 0.3641487431660444, 0.36311204835549, 0.3621421965634748, 0.3614510230013094, 0.360931770810437, 0.3604507824454752, 0.3600321380901917)
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $row = $i + 6
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
}

# --- Column E: updated cf10_simul values (rows 6-33) ---
$eVals = @(0.2247201262549182, 0.3912824475614967, 0.3395398806115439, 0.3257092640157674, 0.3343200114418086, 0.3338690306575793, 0.33462970462644, 0.3379743942214892, 0.3408811792555912, 0.3435454231495588, 0.3459421919281677, 0.3479488585916479, 0.3497068064234699, 0.3512181014165098, 0.3524718136786171, 0.3534973479135821, 0.3543371273321602, 0.3550377706472655, 0.3556134897586724, 0.3560772666050629, 0.3564589122904928, 0.3567751538060301, 0.357033502116559, 0.3572447183045518, 0.3574185979830088, 0.3575620753469036, 0.3576800412029104, 0.3577767497126343)
for ($i = 0; $i -lt $eVals.Length; $i++) {
    $row = $i + 6
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
}
